$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new data row (row 4) with the same shape/style as the existing rows.
$row = 4

# Copy the date-formatted style from A3 so the new date cell reuses the
# existing style record instead of minting a new one.
$ws.Cells.Item(3, 1).Copy()
$ws.Cells.Item($row, 1).PasteSpecial(-4122)
$ws.Cells.Item($row, 1).Value = 42602.583657407406

$ws.Cells.Item($row, 2).Value = "Named"
$ws.Cells.Item($row, 3).Value = 12213
$ws.Cells.Item($row, 4).Value = 6345
$ws.Cells.Item($row, 5).Value = 384
$ws.Cells.Item($row, 6).Value = 85
$ws.Cells.Item($row, 7).Value = 28
$ws.Cells.Item($row, 8).Value = 74
$ws.Cells.Item($row, 9).Value = 24
$ws.Cells.Item($row, 10).Value = 0
$ws.Cells.Item($row, 11).Value = 1
$ws.Cells.Item($row, 12).Value = 0
$ws.Cells.Item($row, 13).Value = 100
